# Update the "clasificacion" precipitation table on Hoja1:
#  - header for column B becomes "Categoría (mm)"
#  - category values lose their "a)/b)/c)/d)/e)" letter prefix and the
#    " mm" suffix (just the numeric range remains, with a trailing space)
#  - the two lowest bands (previously ids 1 and 2, tied on "Baja") swap
#    places so the id numbering matches the new sort order
#  - selection / sort bookkeeping is refreshed to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Categoría (mm)"

# Category (column B) text reformatted - no letter prefix, no "mm" suffix
$ws.Range("B2").Value = "1081-1233 "
$ws.Range("B3").Value = "930-1081 "
$ws.Range("B4").Value = "778-930 "

# Rows 5 and 6 swap which id/category they hold (both stay "Baja")
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "626-778 "
$ws.Range("C5").Value = "Baja"

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "551-626 "
$ws.Range("C6").Value = "Baja"

# Refresh the remembered sort (now only rows 2-5 are considered - the
# lowest, tied "Baja" row 6 sits outside the last sort range)
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("D2:D5"), 0, 2) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:D5")) | Out-Null
$ws.Sort.Header = 0
$ws.Sort.Apply() | Out-Null

# Restore the active cell/selection as last left by the author
$ws.Range("B12").Select() | Out-Null
